$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates ---
$ws.Range("C3").Value = 6.8702204244412304
$ws.Range("D3").Value = [double]"1.12630817165716E-2"
$ws.Range("F3").Value = 3.3821076887715602
$ws.Range("G3").Value = [double]"7.1211759809869296E-2"
$ws.Range("G3").Font.Italic = $true
$ws.Range("I3").Value = 0.227089476136204
$ws.Range("J3").Value = 0.63554460399277801
$ws.Range("L3").Value = 3.9301080217375999
$ws.Range("M3").Value = [double]"5.2724090475718903E-2"
$ws.Range("M3").Font.Italic = $true
$ws.Range("O3").Value = 3.98523787267575
$ws.Range("P3").Value = [double]"5.1148950628178799E-2"
$ws.Range("P3").Font.Italic = $true

# --- Row 4 updates ---
$ws.Range("C4").Value = 2.3358724845536498
$ws.Range("D4").Value = 0.132053362595801
$ws.Range("F4").Value = 3.75408782549344
$ws.Range("G4").Value = [double]"5.7727820627292897E-2"
$ws.Range("I4").Value = 4.38817565473515
$ws.Range("J4").Value = [double]"4.0725095207643801E-2"
$ws.Range("L4").Value = 1.7105331462004001
$ws.Range("M4").Value = 0.196668780008731
$ws.Range("O4").Value = 0.20507402240071801
$ws.Range("P4").Value = 0.65253991458734895

# --- Row 5 updates ---
$ws.Range("C5").Value = 0.12805119993318401
$ws.Range("D5").Value = 0.72180651884157399
$ws.Range("F5").Value = 0.62468001540109397
$ws.Range("G5").Value = 0.432647810764277
$ws.Range("I5").Value = [double]"5.9313977746732896E-3"
$ws.Range("J5").Value = 0.93888566801486295
$ws.Range("L5").Value = 0.11411512985054501
$ws.Range("M5").Value = 0.73686687186911703
$ws.Range("O5").Value = 0.71070731126266795
$ws.Range("P5").Value = 0.403070291208113

# --- Row 6 rebuild: B/E/H/K/N keep a value (now 56/56/56/52/52) with the
#     plain header-row style; C/D/F/G/I/J/L/M (the empty stat-placeholder
#     cells) are removed entirely. Each destination is cleared before the
#     format-only paste so the paste isn't short-circuited by the cell's
#     pre-existing (non-default) style. ---
$ws.Range("B6").Clear()
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4104)
$ws.Range("B6").Value = 56

$ws.Range("E6").Clear()
$ws.Range("B2").Copy()
$ws.Range("E6").PasteSpecial(-4104)
$ws.Range("E6").Value = 56

$ws.Range("H6").Clear()
$ws.Range("B2").Copy()
$ws.Range("H6").PasteSpecial(-4104)
$ws.Range("H6").Value = 56

$ws.Range("K6").Clear()
$ws.Range("B2").Copy()
$ws.Range("K6").PasteSpecial(-4104)
$ws.Range("K6").Value = 52

$ws.Range("N6").Clear()
$ws.Range("B2").Copy()
$ws.Range("N6").PasteSpecial(-4104)
$ws.Range("N6").Value = 52

$ws.Range("C6").Clear()
$ws.Range("D6").Clear()
$ws.Range("F6").Clear()
$ws.Range("G6").Clear()
$ws.Range("I6").Clear()
$ws.Range("J6").Clear()
$ws.Range("L6").Clear()
$ws.Range("M6").Clear()
